$d = $word.ActiveDocument

function Replace-ParagraphRuns($findText, $newInnerXml) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
    $p = $rng.Paragraphs(1)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    # Exclude the trailing paragraph mark from the replaced range.
    $full = $d.Range($pStart, $pEnd - 1)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
           $newInnerXml + `
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $full.InsertXML($xml)
}

# 1. "Strong sense of development pipeline in popular game engines (Unity/SDL/Godot/Unreal Engine)"
#    -> "Good sense of development pipeline in popular game engines (Unity/SDL/Godot/Unreal Engine)"
#    The leading word run is split: "Strong" becomes "Good" (keeping the original run's rsid),
#    and the rest of that original run becomes its own new run.
$inner1 = '<w:r w:rsidRPr="006F0357"><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Good</w:t></w:r>' + `
          '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> sense of development pipeline in popular game engines (Unity/SDL</w:t></w:r>' + `
          '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>/Godot/Unreal Engine</w:t></w:r>' + `
          '<w:r w:rsidRPr="006F0357"><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>)</w:t></w:r>'
Replace-ParagraphRuns "Strong sense of development pipeline in popular game engines (Unity/SDL/Godot/Unreal Engine)" $inner1

# 2. "Good understanding of the software development life cycle"
#    -> "Strong understanding of the software development life cycle"
$inner2 = '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Strong</w:t></w:r>' + `
          '<w:r w:rsidRPr="006F0357"><w:rPr><w:rFonts w:ascii="Century Gothic" w:eastAsia="Century Gothic" w:hAnsi="Century Gothic" w:cs="Century Gothic"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> understanding of the software development life cycle</w:t></w:r>'
Replace-ParagraphRuns "Good understanding of the software development life cycle" $inner2
